$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FeatureSpecifications")
$ws.Columns("Z").Delete()
